# Auto-generated Excel COM-interop edit script
# Updates patient record fields on the active worksheet to match the
# corrected admission/discharge report (v2.0.2: fixed reprint date/time).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A6").Value = "REYES"
$ws.Range("C6").Value = "CASTAÑEDA"
$ws.Range("E6").Value = "SAIDA"
$ws.Range("G6").Value = "LETICIA "
$ws.Range("I6").Value = "2017-20627/201762623"
$ws.Range("A8").Value = ""
$ws.Range("D8").Value = ""
$ws.Range("F8").Value = ""
$ws.Range("H8").Value = ""
$ws.Range("J8").Value = ""
$ws.Range("A10").Value = "ALDEA MAGDALENA "
$ws.Range("D10").Value = ""
$ws.Range("F10").Value = "EL PROGRESO "
$ws.Range("H10").Value = "SAN AGUSTIN "
$ws.Range("J10").Value = ""
$ws.Range("A12").Value = "12 MAYO 1966"
$ws.Range("F12").NumberFormat = "@"
$ws.Range("F12").Value = "51"
$ws.Range("H12").Value = "EL PROGRESO "
$ws.Range("D14").Value = "AMA DE CASA "
$ws.Range("F14").Value = "GUATEMALA"
$ws.Range("H14").NumberFormat = "@"
$ws.Range("H14").Value = "2345336330203"
$ws.Range("A16").Value = "MARCO ANTONIO GIRON "
$ws.Range("F16").Value = "IDEM "
$ws.Range("A18").Value = "REYES HUMBERTO "
$ws.Range("F18").Value = "CASTAÑEDA JOSEFINA "
$ws.Range("A20").Value = "MARCO ANTONIO CHAN "
$ws.Range("F20").Value = "HIJO "
$ws.Range("H20").Value = "IDEM "
$ws.Range("J20").NumberFormat = "@"
$ws.Range("J20").Value = "42900387"
$ws.Range("A24").Value = "24/10/2017"
$ws.Range("C24").Value = "14:41:30"
$ws.Range("D24").Value = "HO"
